$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C7").Copy()
$ws.Range("F7:G7").PasteSpecial(-4122)
$ws.Range("F7").Value = "ок"
$ws.Range("G7").Value = "ок"
